$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the top of the data block (row 69),
# pushing the existing rows 69-75 down to 70-76.
$ws.Rows(69).Insert()

$ws.Cells.Item(69, 1).Value = 3
$ws.Cells.Item(69, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(69, 3).Value = "Coquimbo"
$ws.Cells.Item(69, 4).Value = 45209
$ws.Cells.Item(69, 5).Value = 5
$ws.Cells.Item(69, 6).Value = 300000000
$ws.Cells.Item(69, 7).Value = "Espárragos"
$ws.Cells.Item(69, 8).Value = "Verde"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 1800
$ws.Cells.Item(69, 11).Value = 1800
$ws.Cells.Item(69, 12).Value = 1800
$ws.Cells.Item(69, 13).Value = 1800
$ws.Cells.Item(69, 14).Value = "`$/kilo"
$ws.Cells.Item(69, 15).Value = "Provincia de Linares"
$ws.Cells.Item(69, 16).Value = 1800
$ws.Cells.Item(69, 17).Value = 1
$ws.Cells.Item(69, 18).Value = "Hortaliza"
